$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh: update Price (D) and Volume(1h) (E) columns for rows 2-51.
# Price values that look like plain numbers must stay TEXT (exact digits incl. trailing
# zeros, e.g. "0.7510"), so those cells are pre-formatted as Text before the write —
# otherwise Excel would silently reinterpret them as floating point numbers.

$ws.Range("D2").Value = '26.933.19'
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("D3").Value = '1.863.00'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.94'
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5048'
$ws.Range("E7").Value = '  -1.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3646'
$ws.Range("E8").Value = '  -2.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07176'
$ws.Range("E9").Value = '  +0.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8923'
$ws.Range("E10").Value = '  +0.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.75'
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("D12").Value = '1.862.40'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07483'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.83'
$ws.Range("E14").Value = '  +5.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.231'
$ws.Range("E15").Value = '  -1.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9998'
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008499'
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.21'
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").Value = '26.991.17'
$ws.Range("E20").Value = '  -0.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.029'
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D22").Value = '2.096.87'
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("E23").Value = '  -1.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.417'
$ws.Range("E24").Value = '  -0.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.57'
$ws.Range("E25").Value = '  -1.36%  '
$ws.Range("E26").Value = '  -2.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.88'
$ws.Range("E27").Value = '  -0.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.076'
$ws.Range("E28").Value = '  -0.69%  '
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.703'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.678'
$ws.Range("E31").Value = '  +0.32%  '
$ws.Range("E32").Value = '  +2.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05102'
$ws.Range("E33").Value = '  -0.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7510'
$ws.Range("E34").Value = '  +2.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.998'
$ws.Range("E35").Value = '  -3.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.155'
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.267'
$ws.Range("E37").Value = '  +6.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.540'
$ws.Range("E38").Value = '  +1.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01999'
$ws.Range("E39").Value = '  -2.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5561'
$ws.Range("E40").Value = '  +4.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.075'
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '118.15'
$ws.Range("E42").Value = '  +1.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.546'
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.555'
$ws.Range("E44").Value = '  +2.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1471'
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4692'
$ws.Range("E46").Value = '  +1.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9994'
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.04'
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.565'
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.77'
$ws.Range("E50").Value = '  +0.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.09'
$ws.Range("E51").Value = '  -1.85%  '
